$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 7
$ws.Range("I8").Value = 7
$ws.Range("K8").Value = 21
$ws.Range("M8").Value = 118
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H53").Value = 214.36363
$ws.Range("I53").Value = 236.28572
$ws.Range("K53").Value = 236.28572
$ws.Range("M53").Value = 400.71428
$ws.Range("H55").Value = 1183.1428
$ws.Range("I55").Value = 1801.5
$ws.Range("J55").Value = 358.66666
$ws.Range("K55").Value = 1801.5
$ws.Range("L55").Value = 358.66666
$ws.Range("M55").Value = -1587.5
$ws.Range("N55").Value = -786.66666
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 5000
$ws.Range("K62").Value = 5000
$ws.Range("M62").Value = -4376
$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 5000
$ws.Range("K65").Value = 25000
$ws.Range("M65").Value = -21880

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 750
$ws.Range("I3").Value = 750
$ws.Range("K3").Value = 750
$ws.Range("M3").Value = -635
$ws.Range("H22").Value = 3580.6667
$ws.Range("J22").Value = 5000
$ws.Range("L22").Value = 5000
$ws.Range("N22").Value = -5598
$ws.Range("H41").Value = 7207.7144
$ws.Range("J41").Value = 12124.75
$ws.Range("L41").Value = 12124.75
$ws.Range("N41").Value = -12952.75
$ws.Range("H122").Value = 8899.833000000001
$ws.Range("I122").Value = 8799.799999999999
$ws.Range("J122").Value = 9400
$ws.Range("K122").Value = 26399.4
$ws.Range("L122").Value = 28200
$ws.Range("M122").Value = -23949.4
$ws.Range("N122").Value = -33100
$ws.Range("H132").Value = 3712
$ws.Range("J132").Value = 3392.5
$ws.Range("L132").Value = 10177.5
$ws.Range("N132").Value = -15237.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H25").Value = 1000.5
$ws.Range("I25").Value = 1000.5
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 1000.5
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -765.5
$ws.Range("N25").ClearContents()
$ws.Range("H54").Value = 2943.3333
$ws.Range("I54").Value = 2943.3333
$ws.Range("K54").Value = 2943.3333
$ws.Range("M54").Value = -2459.3333
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H7").Value = 61.75862
$ws.Range("I7").Value = 70.84999999999999
$ws.Range("J7").Value = 41.555557
$ws.Range("K7").Value = 70.84999999999999
$ws.Range("L7").Value = 41.555557
$ws.Range("M7").Value = 42.15000000000001
$ws.Range("N7").Value = -267.555557
$ws.Range("H16").Value = 1122.375
$ws.Range("I16").Value = 696.6667
$ws.Range("J16").Value = 2399.5
$ws.Range("K16").Value = 696.6667
$ws.Range("L16").Value = 2399.5
$ws.Range("M16").Value = -409.6667
$ws.Range("N16").Value = -2973.5
$ws.Range("H22").Value = 1500
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1500
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -2200
$ws.Range("H60").Value = 16611.875
$ws.Range("J60").Value = 25697.5
$ws.Range("L60").Value = 25697.5
$ws.Range("N60").Value = -26719.5
$ws.Range("H86").Value = 20000
$ws.Range("J86").Value = 20000
$ws.Range("L86").Value = 20000
$ws.Range("N86").Value = -22246
$ws.Range("H89").Value = 20000
$ws.Range("J89").Value = 20000
$ws.Range("L89").Value = 100000
$ws.Range("N89").Value = -111232
$ws.Range("H113").Value = 1122.375
$ws.Range("I113").Value = 696.6667
$ws.Range("J113").Value = 2399.5
$ws.Range("K113").Value = 696.6667
$ws.Range("L113").Value = 2399.5
$ws.Range("M113").Value = 1473.3333
$ws.Range("N113").Value = -6739.5
$ws.Range("H132").Value = 7655.8
$ws.Range("I132").Value = 7397
$ws.Range("K132").Value = 22191
$ws.Range("M132").Value = -19661

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3074.4614
$ws.Range("J5").Value = 5496
$ws.Range("L5").Value = 16488
$ws.Range("N5").Value = -16712
$ws.Range("H16").Value = 2500.5
$ws.Range("I16").Value = 2500.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 7501.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -7328.5
$ws.Range("N16").ClearContents()
$ws.Range("H38").Value = 3333.3333
$ws.Range("I38").Value = 2500
$ws.Range("J38").Value = 3750
$ws.Range("K38").Value = 7500
$ws.Range("L38").Value = 11250
$ws.Range("M38").Value = -7153
$ws.Range("N38").Value = -11944
$ws.Range("H50").Value = 500.18182
$ws.Range("I50").Value = 333.83334
$ws.Range("J50").Value = 699.8
$ws.Range("K50").Value = 1001.50002
$ws.Range("L50").Value = 2099.4
$ws.Range("M50").Value = -520.5000200000001
$ws.Range("N50").Value = -3061.4
$ws.Range("H53").Value = 500.18182
$ws.Range("I53").Value = 333.83334
$ws.Range("J53").Value = 699.8
$ws.Range("K53").Value = 1001.50002
$ws.Range("L53").Value = 2099.4
$ws.Range("M53").Value = -520.5000200000001
$ws.Range("N53").Value = -3061.4
$ws.Range("H114").Value = 1873.75
$ws.Range("I114").Value = 2198.3333
$ws.Range("J114").Value = 900
$ws.Range("K114").Value = 6594.999899999999
$ws.Range("L114").Value = 2700
$ws.Range("M114").Value = -3340.999899999999
$ws.Range("N114").Value = -9208
$ws.Range("H117").Value = 7086
$ws.Range("I117").Value = 629
$ws.Range("K117").Value = 1887
$ws.Range("M117").Value = 1555
$ws.Range("H121").Value = 580.4286
$ws.Range("I121").Value = 257.5
$ws.Range("J121").Value = 1011
$ws.Range("K121").Value = 772.5
$ws.Range("L121").Value = 3033
$ws.Range("M121").Value = 537.5
$ws.Range("N121").Value = -5653
$ws.Range("H131").Value = 1013
$ws.Range("J131").Value = 1000
$ws.Range("L131").Value = 3000
$ws.Range("N131").Value = -13080
$ws.Range("H135").Value = 3074.4614
$ws.Range("J135").Value = 5496
$ws.Range("L135").Value = 49464
$ws.Range("N135").Value = -54534
$ws.Range("H139").Value = 46899.09
$ws.Range("I139").Value = 2648.3333
$ws.Range("K139").Value = 7944.999899999999
$ws.Range("M139").Value = -2804.999899999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 5252.25
$ws.Range("I22").Value = 504.5
$ws.Range("J22").Value = 10000
$ws.Range("K22").Value = 504.5
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = 24.5
$ws.Range("N22").Value = -11058
$ws.Range("H41").Value = 11675.333
$ws.Range("I41").Value = 10220.4
$ws.Range("K41").Value = 10220.4
$ws.Range("M41").Value = -9865.4
$ws.Range("H86").Value = 58000
$ws.Range("J86").Value = 58000
$ws.Range("L86").Value = 58000
$ws.Range("N86").Value = -60372
$ws.Range("H89").Value = 58000
$ws.Range("J89").Value = 58000
$ws.Range("L89").Value = 174000
$ws.Range("N89").Value = -185856
$ws.Range("H113").Value = 1999
$ws.Range("I113").Value = 1999
$ws.Range("K113").Value = 1999
$ws.Range("M113").Value = 171
$ws.Range("H122").Value = 5766.727
$ws.Range("I122").Value = 5398.8
$ws.Range("K122").Value = 16196.4
$ws.Range("M122").Value = -13746.4
$ws.Range("H132").Value = 4599.5454
$ws.Range("I132").Value = 4199.5
$ws.Range("K132").Value = 12598.5
$ws.Range("M132").Value = -10068.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H34").Value = 10000
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H68").Value = 3899.75
$ws.Range("I68").Value = 2700
$ws.Range("J68").Value = 5099.5
$ws.Range("K68").Value = 2700
$ws.Range("L68").Value = 5099.5
$ws.Range("M68").Value = -1951
$ws.Range("N68").Value = -6597.5
$ws.Range("H71").Value = 3899.75
$ws.Range("I71").Value = 2700
$ws.Range("J71").Value = 5099.5
$ws.Range("K71").Value = 13500
$ws.Range("L71").Value = 25497.5
$ws.Range("M71").Value = -9756
$ws.Range("N71").Value = -32985.5
$ws.Range("H132").Value = 835125
$ws.Range("I132").Value = 5000000
$ws.Range("K132").Value = 15000000
$ws.Range("M132").Value = -14997470
$ws.Range("H136").Value = 750071.3
$ws.Range("I136").Value = 1020499.8
$ws.Range("K136").Value = 3061499.4
$ws.Range("M136").Value = -3058949.4

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 501000
$ws.Range("I2").Value = 501000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 501000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -500888
$ws.Range("N2").ClearContents()
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H132").Value = 1408
$ws.Range("I132").Value = 1112
$ws.Range("K132").Value = 3336
$ws.Range("M132").Value = -806
